# Updated cryptos list on Thu Nov  2 03:09:37 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns scraped from coinranking.com,
# and reorders two pairs of rows whose ranking flipped (Polygon/WrappedEther at
# rows 14-15, Kaspa/FraxShare at rows 47-48) by swapping their Coin/Link/Price/
# Volume cells in place.
#
# NumberFormat is forced to Text ("@") before writing any Price value that would
# otherwise be auto-parsed as a number by Excel (e.g. "232.12" or "0.0701"),
# matching the worksheet's original plain-text price formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.818.23'
$ws.Range('E2').Value = '  +3.78%  '
$ws.Range('D3').Value = '1.869.16'
$ws.Range('E3').Value = '  +3.10%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.12'
$ws.Range('E5').Value = '  +2.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  +3.00%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.63'
$ws.Range('E8').Value = '  +11.28%  '
$ws.Range('E9').Value = '  +7.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0701'
$ws.Range('E10').Value = '  +3.45%  '
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('D12').Value = '2.140.25'
$ws.Range('E12').Value = '  +3.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.66'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.862.52'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.685'
$ws.Range('E15').Value = '  +8.10%  '
$ws.Range('E16').Value = '  +8.28%  '
$ws.Range('D17').Value = '35.827.94'
$ws.Range('E17').Value = '  +3.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.56'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').Value = '0.0₃0807'
$ws.Range('E19').Value = '  +4.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '247.09'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.45'
$ws.Range('E21').Value = '  +10.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.79'
$ws.Range('E22').Value = '  +15.84%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.24'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.93'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.06'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.46'
$ws.Range('E29').Value = '  +18.54%  '
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '3.342.89'
$ws.Range('E31').Value = '  +37.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0549'
$ws.Range('E32').Value = '  +5.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.97'
$ws.Range('E33').Value = '  +4.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.06'
$ws.Range('E34').Value = '  +4.89%  '
$ws.Range('E35').Value = '  +4.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '97.99'
$ws.Range('E36').Value = '  +19.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.691'
$ws.Range('E37').Value = '  +6.63%  '
$ws.Range('E38').Value = '  +7.23%  '
$ws.Range('D39').Value = '1.357.46'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0198'
$ws.Range('E41').Value = '  +5.78%  '
$ws.Range('E42').Value = '  +7.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.25'
$ws.Range('E43').Value = '  +10.49%  '
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.24'
$ws.Range('E47').Value = '  +7.83%  '
$ws.Range('B48').Value = 'Kaspa'
$ws.Range('C48').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0520'
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').Value = '2.037.45'
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '105.28'
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('E51').Value = '  +0.24%  '
